$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4121
$ws.Range("I62").Value = 3377
$ws.Range("K62").Value = 3377
$ws.Range("M62").Value = -2753
$ws.Range("H64").Value = 400002700
$ws.Range("I64").Value = 4500
$ws.Range("K64").Value = 4500
$ws.Range("M64").Value = -4252
$ws.Range("H65").Value = 4121
$ws.Range("I65").Value = 3377
$ws.Range("K65").Value = 16885
$ws.Range("M65").Value = -13765
$ws.Range("H67").Value = 400002700
$ws.Range("I67").Value = 4500
$ws.Range("K67").Value = 4500
$ws.Range("M67").Value = -3642
$ws.Range("H132").Value = 6446.5654
$ws.Range("I132").Value = 6970.5713
$ws.Range("K132").Value = 20911.7139
$ws.Range("M132").Value = -18381.7139
$ws.Range("H137").Value = 2916.05
$ws.Range("I137").Value = 1939.6364
$ws.Range("K137").Value = 5818.9092
$ws.Range("M137").Value = -3268.9092
$ws.Range("H138").Value = 6466.4443
$ws.Range("J138").Value = 7199.6665
$ws.Range("L138").Value = 21598.9995
$ws.Range("N138").Value = -31878.9995
$ws.Range("H141").Value = 2997.8667
$ws.Range("I141").Value = 2533.182
$ws.Range("J141").Value = 4275.75
$ws.Range("K141").Value = 7599.545999999999
$ws.Range("L141").Value = 12827.25
$ws.Range("M141").Value = -2419.545999999999
$ws.Range("N141").Value = -23187.25

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2264.182
$ws.Range("I2").Value = 2099.2222
$ws.Range("K2").Value = 2099.2222
$ws.Range("M2").Value = -1986.2222
$ws.Range("H32").Value = 1389276.5
$ws.Range("I32").Value = 651295.4
$ws.Range("J32").Value = 11905506
$ws.Range("K32").Value = 651295.4
$ws.Range("L32").Value = 11905506
$ws.Range("M32").Value = -651008.4
$ws.Range("N32").Value = -11906080
$ws.Range("H110").Value = 669.75
$ws.Range("I110").Value = 405
$ws.Range("K110").Value = 405
$ws.Range("M110").Value = 1640
$ws.Range("H116").Value = 2264.182
$ws.Range("I116").Value = 2099.2222
$ws.Range("K116").Value = 2099.2222
$ws.Range("M116").Value = 194.7777999999998
$ws.Range("H122").Value = 3668.0908
$ws.Range("I122").Value = 3549.8
$ws.Range("J122").Value = 3766.6667
$ws.Range("K122").Value = 10649.4
$ws.Range("L122").Value = 11300.0001
$ws.Range("M122").Value = -8199.400000000001
$ws.Range("N122").Value = -16200.0001
$ws.Range("H132").Value = 1396.9286
$ws.Range("I132").Value = 1090.4138
$ws.Range("K132").Value = 3271.2414
$ws.Range("M132").Value = -741.2413999999999

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2264.182
$ws.Range("I3").Value = 2099.2222
$ws.Range("K3").Value = 2099.2222
$ws.Range("M3").Value = -1985.2222
$ws.Range("H132").Value = 109992.5
$ws.Range("J132").Value = 109992.5
$ws.Range("L132").Value = 109992.5
$ws.Range("N132").Value = -120112.5
$ws.Range("H134").Value = 3237.25
$ws.Range("I134").Value = 2900
$ws.Range("K134").Value = 8700
$ws.Range("M134").Value = -6165

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11372364
$ws.Range("I31").Value = 2668
$ws.Range("K31").Value = 2668
$ws.Range("M31").Value = -2373
$ws.Range("H34").Value = 11372364
$ws.Range("I34").Value = 2668
$ws.Range("K34").Value = 2668
$ws.Range("M34").Value = -2466
$ws.Range("H99").Value = 2840.16
$ws.Range("I99").Value = 2139.0908
$ws.Range("J99").Value = 3391
$ws.Range("K99").Value = 2139.0908
$ws.Range("L99").Value = 3391
$ws.Range("M99").Value = -641.0907999999999
$ws.Range("N99").Value = -6387
$ws.Range("H122").Value = 2636.5264
$ws.Range("J122").Value = 3122.111
$ws.Range("L122").Value = 9366.332999999999
$ws.Range("N122").Value = -14266.333
$ws.Range("H126").Value = 2840.16
$ws.Range("I126").Value = 2139.0908
$ws.Range("J126").Value = 3391
$ws.Range("K126").Value = 6417.2724
$ws.Range("L126").Value = 10173
$ws.Range("M126").Value = -3947.2724
$ws.Range("N126").Value = -15113
$ws.Range("H132").Value = 3709.05
$ws.Range("I132").Value = 2772.5
$ws.Range("K132").Value = 8317.5
$ws.Range("M132").Value = -5787.5
$ws.Range("H134").Value = 3421.2942
$ws.Range("I134").Value = 3344.2334
$ws.Range("K134").Value = 10032.7002
$ws.Range("M134").Value = -7497.700199999999
$ws.Range("H141").Value = 165000
$ws.Range("J141").Value = 165000
$ws.Range("L141").Value = 165000
$ws.Range("N141").Value = -175360

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 4302.7144
$ws.Range("I114").Value = 3379.75
$ws.Range("J114").Value = 5533.3335
$ws.Range("K114").Value = 10139.25
$ws.Range("L114").Value = 16600.0005
$ws.Range("M114").Value = -6885.25
$ws.Range("N114").Value = -23108.0005
$ws.Range("H121").Value = 14457536
$ws.Range("I121").Value = 33366916
$ws.Range("J121").Value = 275499.75
$ws.Range("K121").Value = 100100748
$ws.Range("L121").Value = 826499.25
$ws.Range("M121").Value = -100099438
$ws.Range("N121").Value = -829119.25
$ws.Range("H122").Value = 552
$ws.Range("I122").Value = 551.6667
$ws.Range("K122").Value = 4965.0003
$ws.Range("M122").Value = -2515.0003
$ws.Range("H134").Value = 3604.3333
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3500107.5
$ws.Range("I122").Value = 5497310.5
$ws.Range("K122").Value = 16491931.5
$ws.Range("M122").Value = -16489481.5
$ws.Range("H132").Value = 2129.9143
$ws.Range("I132").Value = 2045.2609
$ws.Range("K132").Value = 6135.7827
$ws.Range("M132").Value = -3605.7827

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 788.0909
$ws.Range("I16").Value = 788.0909
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 788.0909
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -618.0909
$ws.Range("N16").ClearContents()
$ws.Range("H122").Value = 3257.1875
$ws.Range("I122").Value = 2523.25
$ws.Range("J122").Value = 3501.8333
$ws.Range("K122").Value = 7569.75
$ws.Range("L122").Value = 10505.4999
$ws.Range("M122").Value = -5119.75
$ws.Range("N122").Value = -15405.4999

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 12501905
$ws.Range("I122").Value = 1822
$ws.Range("J122").Value = 35716344
$ws.Range("K122").Value = 5466
$ws.Range("L122").Value = 107149032
$ws.Range("M122").Value = -3016
$ws.Range("N122").Value = -107153932
